$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ikisi de değil"
$ws.Range("B1").Value = "180 TL"
$ws.Range("C1").Value = "290 TL"

$ws.Range("D1").Value = "'"
$ws.Range("D1").Style = "Normal"

$ws.Range("E1").Value = "['https://cdn.dsmcdn.com/mnresize/1200/1800/ty1157/product/media/images/prod/SPM/PIM/20240202/12/61fa1b61-7552-3f42-9812-6eac2c073ca1/1_org_zoom.jpg', 'https://cdn.dsmcdn.com/ty1157/product/media/images/prod/SPM/PIM/20240202/12/61fa1b61-7552-3f42-9812-6eac2c073ca1/1_org_zoom.jpg']"

$ws.Range("F1").Value = "'"
$ws.Range("F1").Style = "Normal"

$ws.Range("G1").Value = "'"
$ws.Range("G1").Style = "Normal"
